# ----------------------------------------------------------------------
# Restructure "Input" sheet (sheet 1) to the new standard template layout
# Old columns: A=거래처명 B=현장명 C=발주일 D=납기일 E=발주번호 F=품목 G=규격
#              H=수량 I=단위 J=단가 K=공급가액 L=부가세 M=합계 N=대분류
#              O=중분류 P=소분류 Q=비고
# New columns: A=발주일자 B=납기일자 C=거래처명 D=거래처이메일 E=납품처명
#              F=납품처이메일 G=프로젝트명 H=대분류 I=중분류 J=소분류 K=품목명
#              L=규격 M=수량 N=단가 O=총금액 P=비고
# ----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    # --- read old values (by old column letter) ---
    # NOTE: this runtime's `.Value` GETTER returns a bogus reflection
    #       descriptor string instead of the cell contents (the SETTER is
    #       fine) - so reads must go through `.Value2` instead.
    $oldA = $ws.Cells.Item($r, 1).Value2   # 거래처명
    $oldB = $ws.Cells.Item($r, 2).Value2   # 현장명
    $oldC = $ws.Cells.Item($r, 3).Value2   # 발주일
    $oldD = $ws.Cells.Item($r, 4).Value2   # 납기일
    $oldF = $ws.Cells.Item($r, 6).Value2   # 품목
    $oldG = $ws.Cells.Item($r, 7).Value2   # 규격
    $oldH = $ws.Cells.Item($r, 8).Value2   # 수량
    $oldJ = $ws.Cells.Item($r, 10).Value2  # 단가
    $oldM = $ws.Cells.Item($r, 13).Value2  # 합계
    $oldN = $ws.Cells.Item($r, 14).Value2  # 대분류
    $oldO = $ws.Cells.Item($r, 15).Value2  # 중분류
    $oldP = $ws.Cells.Item($r, 16).Value2  # 소분류

    # --- derived new fields ---
    $buyerEmail = "$oldA@example.com"
    $deliveryEmail = "delivery@example.com"

    # --- clear the whole row first: guarantees no stale leftover values
    #     (e.g. old P/Q) survive in columns that have no new source, and
    #     strips any pre-existing cell style (data rows carry none in the
    #     target layout). Assigning "" (rather than ClearContents()) makes
    #     the exporter drop the <c> element entirely instead of leaving an
    #     empty placeholder behind. ---
    for ($cc = 1; $cc -le 17; $cc++) {
        $ws.Cells.Item($r, $cc).Value = ""
    }

    # --- write new layout. 발주일자/납기일자 are free-text date-looking
    #     strings ("2025-08-21") - prefixing with a single-quote stops
    #     Excel's auto-detection from silently turning them into date
    #     serial numbers; ClearFormats() afterwards drops the resulting
    #     quote-prefix cell style again so the cell ends up unstyled. ---
    $ws.Cells.Item($r, 1).Value = "'" + $oldC   # A: 발주일자   (was C 발주일)
    $ws.Cells.Item($r, 2).Value = "'" + $oldD   # B: 납기일자   (was D 납기일)
    $ws.Cells.Item($r, 3).Value = $oldA          # C: 거래처명   (was A 거래처명)
    $ws.Cells.Item($r, 4).Value = $buyerEmail    # D: 거래처 이메일 (new, derived)
    $ws.Cells.Item($r, 5).Value = $oldB          # E: 납품처명   (was B 현장명)
    $ws.Cells.Item($r, 6).Value = $deliveryEmail # F: 납품처 이메일 (new constant)
    $ws.Cells.Item($r, 7).Value = $oldB          # G: 프로젝트명 (was B 현장명)
    $ws.Cells.Item($r, 8).Value = $oldN          # H: 대분류     (was N 대분류)
    $ws.Cells.Item($r, 9).Value = $oldO          # I: 중분류     (was O 중분류)
    $ws.Cells.Item($r, 10).Value = $oldP         # J: 소분류     (was P 소분류)
    $ws.Cells.Item($r, 11).Value = $oldF         # K: 품목명     (was F 품목)
    $ws.Cells.Item($r, 12).Value = $oldG         # L: 규격       (was G 규격)
    $ws.Cells.Item($r, 13).Value = $oldH         # M: 수량       (was H 수량, numeric)
    $ws.Cells.Item($r, 14).Value = $oldJ         # N: 단가       (was J 단가, numeric)
    $ws.Cells.Item($r, 15).Value = $oldM         # O: 총금액     (was M 합계, numeric)
    # P (비고) has no source column - stays blank.

    $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,17)).ClearFormats()
}

# --- header row (row 1): new labels, and drop the bold/bordered header style ---
$ws.Cells.Item(1, 1).Value = "발주일자"
$ws.Cells.Item(1, 2).Value = "납기일자"
$ws.Cells.Item(1, 3).Value = "거래처명"
$ws.Cells.Item(1, 4).Value = "거래처 이메일"
$ws.Cells.Item(1, 5).Value = "납품처명"
$ws.Cells.Item(1, 6).Value = "납품처 이메일"
$ws.Cells.Item(1, 7).Value = "프로젝트명"
$ws.Cells.Item(1, 8).Value = "대분류"
$ws.Cells.Item(1, 9).Value = "중분류"
$ws.Cells.Item(1, 10).Value = "소분류"
$ws.Cells.Item(1, 11).Value = "품목명"
$ws.Cells.Item(1, 12).Value = "규격"
$ws.Cells.Item(1, 13).Value = "수량"
$ws.Cells.Item(1, 14).Value = "단가"
$ws.Cells.Item(1, 15).Value = "총금액"
$ws.Cells.Item(1, 16).Value = "비고"
$ws.Cells.Item(1, 17).Value = ""
$ws.Range("A1:Q1").ClearFormats()

# --- drop now-empty column Q entirely so dimension becomes A1:P7 ---
$ws.Columns.Item(17).Delete()

# ----------------------------------------------------------------------
# Sheet 2 ("갑지") and Sheet 3 ("을지"): the trailing remarks column (I)
# was only ever populated with empty strings - clear it out entirely so
# no stray empty cell is serialized.
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
for ($r = 2; $r -le 7; $r++) { $ws2.Cells.Item($r, 9).Value = "" }

$ws3 = $wb.Worksheets.Item(3)
for ($r = 2; $r -le 7; $r++) { $ws3.Cells.Item($r, 9).Value = "" }

Write-Host "Restructure complete"
